$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8999976118043378
$ws.Range("C2").Value = 0.2771270230097684
$ws.Range("E2").Value = 0.1163700776058363
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.3353931556889975
$ws.Range("H2").Value = 0.5015183666473888
$ws.Range("L2").Value = 0.1973744587832158
$ws.Range("N2").Value = 1.046251808964243
$ws.Range("O2").Value = 1.607580560632726

$ws.Range("B3").Value = 0.8086648358936372
$ws.Range("C3").Value = 0.272106846402167
$ws.Range("E3").Value = 0.1170412136229899
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.3351425347788322
$ws.Range("H3").Value = 0.5051237346264017
$ws.Range("L3").Value = 0.1882302716911397
$ws.Range("N3").Value = 1.047611093682185
$ws.Range("O3").Value = 1.614339806777366

$ws.Range("B4").Value = 0.7525909376618642
$ws.Range("C4").Value = 0.2690385543669578
$ws.Range("E4").Value = 0.1175197626729005
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.335293545927648
$ws.Range("H4").Value = 0.5076050336218785
$ws.Range("L4").Value = 0.1827077027227659
$ws.Range("N4").Value = 1.048863099074936
$ws.Range("O4").Value = 1.619687510133119

$ws.Range("B5").Value = 0.7297432065051339
$ws.Range("C5").Value = 0.2677919151960992
$ws.Range("E5").Value = 0.1177314947773667
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3354315907437595
$ws.Range("H5").Value = 0.5086834744196906
$ws.Range("L5").Value = 0.1804804610007977
$ws.Range("N5").Value = 1.0494785267271
$ws.Range("O5").Value = 1.622167523658661

$ws.Range("B6").Value = 0.7259495776141307
$ws.Range("C6").Value = 0.2675851409330932
$ws.Range("E6").Value = 0.1177676627049582
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.3354591290488429
$ws.Range("H6").Value = 0.5088666130449155
$ws.Range("L6").Value = 0.1801120368709377
$ws.Range("N6").Value = 1.049587081300167
$ws.Range("O6").Value = 1.622597484477922

$ws.Range("B7").Value = 0.7522827915157677
$ws.Range("C7").Value = 0.2690217264952963
$ws.Range("E7").Value = 0.1175225504621658
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.3352950980907536
$ws.Range("H7").Value = 0.5076193053866191
$ws.Range("L7").Value = 0.1826775710695756
$ws.Range("N7").Value = 1.048870972537472
$ws.Range("O7").Value = 1.619719739064465

$ws.Range("B8").Value = 0.8685060253064876
$ws.Range("C8").Value = 0.2753932333250049
$ws.Range("E8").Value = 0.1165876927907359
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.3352433628332179
$ws.Range("H8").Value = 0.5027059616239029
$ws.Range("L8").Value = 0.1942025151694793
$ws.Range("N8").Value = 1.046634015979038
$ws.Range("O8").Value = 1.609662400426885

$ws.Range("B9").Value = 1.096396883263594
$ws.Range("C9").Value = 0.2879931478376108
$ws.Range("E9").Value = 0.115281679248735
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.3375689767277379
$ws.Range("H9").Value = 0.4951941392221499
$ws.Range("L9").Value = 0.2175290847482927
$ws.Range("N9").Value = 1.045548529081827
$ws.Range("O9").Value = 1.599458652273682

$ws.Range("B10").Value = 1.263749435957436
$ws.Range("C10").Value = 0.2973066444113783
$ws.Range("E10").Value = 0.114643465033879
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.3407689632092854
$ws.Range("H10").Value = 0.490969943602579
$ws.Range("L10").Value = 0.2351066880623591
$ws.Range("N10").Value = 1.046750640273089
$ws.Range("O10").Value = 1.597790854824837

$ws.Range("B11").Value = 1.339852518879582
$ws.Range("C11").Value = 0.3015542146575427
$ws.Range("E11").Value = 0.1144228803297125
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.3425511164970914
$ws.Range("H11").Value = 0.4893294942384188
$ws.Range("L11").Value = 0.2431980942630787
$ws.Range("N11").Value = 1.047729159130057
$ws.Range("O11").Value = 1.598303556878051

$ws.Range("B12").Value = 1.368665532576074
$ws.Range("C12").Value = 0.3031640595945646
$ws.Range("E12").Value = 0.1143493765307362
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.3432731043544806
$ws.Range("H12").Value = 0.4887487351601862
$ws.Range("L12").Value = 0.2462757072233188
$ws.Range("N12").Value = 1.048161546804394
$ws.Range("O12").Value = 1.598680945541332

$ws.Range("B13").Value = 1.362460410535959
$ws.Range("C13").Value = 0.3028172919501344
$ws.Range("E13").Value = 0.1143647609756115
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.3431155126000505
$ws.Range("H13").Value = 0.488872013073717
$ws.Range("L13").Value = 0.2456122865551578
$ws.Range("N13").Value = 1.048065677477624
$ws.Range("O13").Value = 1.5985915119428

$ws.Range("B14").Value = 1.342223104666004
$ws.Range("C14").Value = 0.3016866308536237
$ws.Range("E14").Value = 0.1144166321869129
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.3426095693897082
$ws.Range("H14").Value = 0.4892809042161446
$ws.Range("L14").Value = 0.2434510200337314
$ws.Range("N14").Value = 1.047763493784103
$ws.Range("O14").Value = 1.598330930112866

$ws.Range("B15").Value = 1.329826392309542
$ws.Range("C15").Value = 0.3009942431849879
$ws.Range("E15").Value = 0.1144497105311082
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.3423058069597147
$ws.Range("H15").Value = 0.4895366291755892
$ws.Range("L15").Value = 0.2421289457431612
$ws.Range("N15").Value = 1.047586444685336
$ws.Range("O15").Value = 1.598195191669276

$ws.Range("B16").Value = 1.258775237138991
$ws.Range("C16").Value = 0.2970292585249439
$ws.Range("E16").Value = 0.1146592839854605
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.3406590800045421
$ws.Range("H16").Value = 0.4910828095880362
$ws.Range("L16").Value = 0.2345798018660048
$ws.Range("N16").Value = 1.046695358338113
$ws.Range("O16").Value = 1.597782967739903

$ws.Range("B17").Value = 1.215179623503673
$ws.Range("C17").Value = 0.2945995153950776
$ws.Range("E17").Value = 0.1148057117653885
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.339732610050163
$ws.Range("H17").Value = 0.4921033620680504
$ws.Range("L17").Value = 0.2299729569779032
$ws.Range("N17").Value = 1.046259102164825
$ws.Range("O17").Value = 1.59785601193704

$ws.Range("B18").Value = 1.190102197293811
$ws.Range("C18").Value = 0.2932030195462971
$ws.Range("E18").Value = 0.1148964977394122
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.3392304466477469
$ws.Range("H18").Value = 0.4927168200343317
$ws.Range("L18").Value = 0.2273321961267385
$ws.Range("N18").Value = 1.046048823665544
$ws.Range("O18").Value = 1.598017669817125

$ws.Range("B19").Value = 1.181611056529107
$ws.Range("C19").Value = 0.2927303719085188
$ws.Range("E19").Value = 0.1149283638433083
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.33906569301692
$ws.Range("H19").Value = 0.4929290710311136
$ws.Range("L19").Value = 0.2264396250891281
$ws.Range("N19").Value = 1.045984614802919
$ws.Range("O19").Value = 1.598092940807589

$ws.Range("B20").Value = 1.219820710397244
$ws.Range("C20").Value = 0.294858060543433
$ws.Range("E20").Value = 0.1147894448823266
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.3398280538920062
$ws.Range("H20").Value = 0.4919919835784015
$ws.Range("L20").Value = 0.2304624355685121
$ws.Range("N20").Value = 1.046301337506208
$ws.Range("O20").Value = 1.597835850452668

$ws.Range("B21").Value = 1.348167457285399
$ws.Range("C21").Value = 0.3020186972835717
$ws.Range("E21").Value = 0.1144011242456919
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.3427568967862129
$ws.Range("H21").Value = 0.4891597052459247
$ws.Range("L21").Value = 0.2440854690402148
$ws.Range("N21").Value = 1.047850575847733
$ws.Range("O21").Value = 1.598402492955245

$ws.Range("B22").Value = 1.432016177386515
$ws.Range("C22").Value = 0.3067065646454807
$ws.Range("E22").Value = 0.1142057773570357
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.3449458284580942
$ws.Range("H22").Value = 0.4875443893371312
$ws.Range("L22").Value = 0.2530679641574096
$ws.Range("N22").Value = 1.049223462214229
$ws.Range("O22").Value = 1.599841076847071

$ws.Range("B23").Value = 1.387268168636865
$ws.Range("C23").Value = 0.3042038885604939
$ws.Range("E23").Value = 0.1143046908097425
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.3437523544542813
$ws.Range("H23").Value = 0.4883849391028434
$ws.Range("L23").Value = 0.2482666456890144
$ws.Range("N23").Value = 1.048457827939117
$ws.Range("O23").Value = 1.598975396376346

$ws.Range("B24").Value = 1.217722517635536
$ws.Range("C24").Value = 0.2947411710365913
$ws.Range("E24").Value = 0.1147967785736128
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.339784808813036
$ws.Range("H24").Value = 0.4920422545361447
$ws.Range("L24").Value = 0.230241118069074
$ws.Range("N24").Value = 1.046282116639063
$ws.Range("O24").Value = 1.597844592747009

$ws.Range("B25").Value = 1.034755526531228
$ws.Range("C25").Value = 0.2845740741551737
$ws.Range("E25").Value = 0.1155785564558638
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.3366787826482849
$ws.Range("H25").Value = 0.4969989424160985
$ws.Range("L25").Value = 0.2111412085871507
$ws.Range("N25").Value = 1.048065677477624
$ws.Range("O25").Value = 1.5985915119428
